$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9300694554254023
$ws.Range("J2").Value = 0.9300694554254023
$ws.Range("M2").Value = 1.321358333333333
$ws.Range("N2").Value = 3.964075
$ws.Range("O2").Value = 0.06904315418552966
$ws.Range("P2").Value = 0.06904315418552966
$ws.Range("Q2").Value = 3.313466345644444
$ws.Range("R2").Value = 29.8211971108
$ws.Range("S2").Value = 0.06421492881418765
$ws.Range("T2").Value = 0.06421492881418765

# Row 3
$ws.Range("I3").Value = 0.9300694554254023
$ws.Range("J3").Value = 0.9300694554254023
$ws.Range("O3").Value = 0.4558096119837698
$ws.Range("P3").Value = 0.4558096119837698
$ws.Range("S3").Value = 0.4239345975954087
$ws.Range("T3").Value = 0.4239345975954087

# Row 4
$ws.Range("I4").Value = 0.9300694554254023
$ws.Range("J4").Value = 0.9300694554254023
$ws.Range("M4").Value = 9.093439666666667
$ws.Range("N4").Value = 27.280319
$ws.Range("O4").Value = 0.4751472338307006
$ws.Range("P4").Value = 0.4751472338307005
$ws.Range("Q4").Value = 22.80290330151289
$ws.Range("R4").Value = 205.226129713616
$ws.Range("S4").Value = 0.441919929015806
$ws.Range("T4").Value = 0.4419199290158059

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1885443333333333
$ws.Range("H5").Value = 0.5656329999999999
$ws.Range("I5").Value = 0.06993054457459773
$ws.Range("J5").Value = 0.06993054457459771
$ws.Range("M5").Value = 1.321358333333333
$ws.Range("N5").Value = 3.964075
$ws.Range("O5").Value = 0.06904315418552966
$ws.Range("P5").Value = 0.06904315418552966
$ws.Range("Q5").Value = 0.2491346260527777
$ws.Range("R5").Value = 2.242211634475
$ws.Range("S5").Value = 0.004828225371342005
$ws.Range("T5").Value = 0.004828225371342004

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1885443333333333
$ws.Range("H6").Value = 0.5656329999999999
$ws.Range("I6").Value = 0.06993054457459773
$ws.Range("J6").Value = 0.06993054457459771
$ws.Range("O6").Value = 0.4558096119837698
$ws.Range("P6").Value = 0.4558096119837698
$ws.Range("Q6").Value = 1.644738838664444
$ws.Range("R6").Value = 14.80264954798
$ws.Range("S6").Value = 0.03187501438836111
$ws.Range("T6").Value = 0.0318750143883611

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.1885443333333333
$ws.Range("H7").Value = 0.5656329999999999
$ws.Range("I7").Value = 0.06993054457459773
$ws.Range("J7").Value = 0.06993054457459771
$ws.Range("M7").Value = 9.093439666666667
$ws.Range("N7").Value = 27.280319
$ws.Range("O7").Value = 0.4751472338307006
$ws.Range("P7").Value = 0.4751472338307005
$ws.Range("Q7").Value = 1.714516519658555
$ws.Range("R7").Value = 15.430648676927
$ws.Range("S7").Value = 0.03322730481489462
$ws.Range("T7").Value = 0.0332273048148946
